$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 14 ("Programa resumido:" label + short-syllabus text in B/C).
# This shifts every row below it up by one, matching the target layout.
$ws.Rows.Item(14).Delete()

# --- Value fixups on top of the shifted layout ---

# Row 10 (Objetivos value) now holds the "Docentes responsaveis" value.
$ws.Range("B10").Value = "5817344 - Livia Melo Carneiro"
$ws.Range("C10").Value = "5817344 - Livia Melo Carneiro"

# Row 13 previously had only B/C ("5817344 - Livia Melo Carneiro").
# It now becomes the "Programa resumido:" row with a "Semestral" value,
# and needs the 60pt custom row height like its neighbours.
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 15 ("Programa:") value becomes "01/01/2012".
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"

# Row 18 ("Metodo:") value becomes the Livia Melo Carneiro string.
$ws.Range("B18").Value = "5817344 - Livia Melo Carneiro"
$ws.Range("C18").Value = "5817344 - Livia Melo Carneiro"

# Row 19 ("Criterio:") value becomes the P1/P2 evaluation text.
$ws.Range("B19").Value = "O aluno será avaliado através de duas provas escritas P1 e P2."
$ws.Range("C19").Value = "O aluno será avaliado através de duas provas escritas P1 e P2."

# Row 20 ("Norma de recuperacao:") value becomes the NF formula text.
$ws.Range("B20").Value = "A nota final NF será (P1 + P2)/2 ."
$ws.Range("C20").Value = "A nota final NF será (P1 + P2)/2 ."

# Row 21 ("Bibliografia:") value becomes the recovery-exam text.
$ws.Range("B21").Value = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
$ws.Range("C21").Value = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
